# Monthly Update Real Estate
# Updates the Zillow Summary Table with the latest month's data and
# moves the active selection from J12 to J11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 1 ---
$ws.Range("A1").Value = 43999

# --- Row 6 ---
$ws.Range("C6").Value = 43891

# --- Row 9 ---
$ws.Range("C9").Value = 658017
$ws.Range("D9").Value = 0.00586382056342205
$ws.Range("E9").Value = 0.0254341619707774
$ws.Range("F9").Value = 0.84560412870327795

# --- Row 10 ---
$ws.Range("C10").Value = 578267
$ws.Range("D10").Value = 0.0074776775991984802
$ws.Range("E10").Value = 0.043929500499157802
$ws.Range("F10").Value = 1.0914301892134199

# --- Row 11 ---
$ws.Range("C11").Value = 248857
$ws.Range("D11").Value = 0.00442361792211043
$ws.Range("E11").Value = 0.040554779685396197
$ws.Range("F11").Value = 1.1628614418422101

# --- Row 14 ---
$ws.Range("C14").Value = 903246
$ws.Range("D14").Value = 0.0062228460664393
$ws.Range("E14").Value = 0.011137393274577301
$ws.Range("F14").Value = 0.93132162160490095

# --- Row 15 ---
$ws.Range("C15").Value = 856454
$ws.Range("D15").Value = 0.0069828280512396201
$ws.Range("E15").Value = 0.030349110411476098
$ws.Range("F15").Value = 0.836244479660055

# --- Row 16 ---
$ws.Range("C16").Value = 368496
$ws.Range("D16").Value = 0.0065914014035068904
$ws.Range("E16").Value = 0.030936833800547101
$ws.Range("F16").Value = 0.73414791099281396

# --- Row 17 ---
$ws.Range("C17").Value = 3276825
$ws.Range("D17").Value = 0.0029984429278100402
$ws.Range("E17").Value = 0.0068745828790837101
$ws.Range("F17").Value = 1.39801298340043

# --- Row 18 ---
$ws.Range("C18").Value = 1134649
$ws.Range("D18").Value = 0.0062673934135166797
$ws.Range("E18").Value = 0.019905546407024102
$ws.Range("F18").Value = 0.95263878615584996

# --- Row 19 ---
$ws.Range("C19").Value = 412599
$ws.Range("D19").Value = 0.0058532708594385401
$ws.Range("E19").Value = 0.039829130911427799
$ws.Range("F19").Value = 0.81380473372781104

# --- Row 20 ---
$ws.Range("C20").Value = 813113
$ws.Range("D20").Value = 0.00162479643850544
$ws.Range("E20").Value = 0.0017025344603076599
$ws.Range("F20").Value = 0.85157613949245603

# --- Active selection moves from J12 to J11 ---
$ws.Range("J11").Select()
